$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntn1"
$ws.Range("C2").Value = "Unc5a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.9305633333333333
$ws.Range("H2").Value = 2.79169
$ws.Range("I2").Value = 0.01768777137856805
$ws.Range("J2").Value = 0.01768777137856806
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.292844666666667
$ws.Range("N2").Value = 6.878534
$ws.Range("O2").Value = 0.4949693416994264
$ws.Range("P2").Value = 0.4949693416994265
$ws.Range("Q2").Value = 2.133637175828889
$ws.Range("R2").Value = 19.20273458246
$ws.Range("S2").Value = 0.008754904555379785
$ws.Range("T2").Value = 0.008754904555379788

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntn1"
$ws.Range("C3").Value = "Unc5a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.9305633333333333
$ws.Range("H3").Value = 2.79169
$ws.Range("I3").Value = 0.01768777137856805
$ws.Range("J3").Value = 0.01768777137856806
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.481271666666667
$ws.Range("N3").Value = 4.443815
$ws.Range("O3").Value = 0.3197704896398035
$ws.Range("P3").Value = 0.3197704896398036
$ws.Range("Q3").Value = 1.378417099705555
$ws.Range("R3").Value = 12.40575389735
$ws.Range("S3").Value = 0.005656027314361609
$ws.Range("T3").Value = 0.005656027314361611

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntn1"
$ws.Range("C4").Value = "Unc5a"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.9305633333333333
$ws.Range("H4").Value = 2.79169
$ws.Range("I4").Value = 0.01768777137856805
$ws.Range("J4").Value = 0.01768777137856806
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8581799999999999
$ws.Range("N4").Value = 2.57454
$ws.Range("O4").Value = 0.18526016866077
$ws.Range("P4").Value = 0.1852601686607701
$ws.Range("Q4").Value = 0.7985908413999999
$ws.Range("R4").Value = 7.1873175726
$ws.Range("S4").Value = 0.003276839508826658
$ws.Range("T4").Value = 0.003276839508826659

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntn1"
$ws.Range("C5").Value = "Unc5a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 44.154177
$ws.Range("H5").Value = 132.462531
$ws.Range("I5").Value = 0.8392647337471152
$ws.Range("J5").Value = 0.8392647337471153
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.292844666666667
$ws.Range("N5").Value = 6.878534
$ws.Range("O5").Value = 0.4949693416994264
$ws.Range("P5").Value = 0.4949693416994265
$ws.Range("Q5").Value = 101.238669245506
$ws.Range("R5").Value = 911.1480232095541
$ws.Range("S5").Value = 0.415410312774354
$ws.Range("T5").Value = 0.4154103127743541

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntn1"
$ws.Range("C6").Value = "Unc5a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 44.154177
$ws.Range("H6").Value = 132.462531
$ws.Range("I6").Value = 0.8392647337471152
$ws.Range("J6").Value = 0.8392647337471153
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.481271666666667
$ws.Range("N6").Value = 4.443815
$ws.Range("O6").Value = 0.3197704896398035
$ws.Range("P6").Value = 0.3197704896398036
$ws.Range("Q6").Value = 65.404331355085
$ws.Range("R6").Value = 588.638982195765
$ws.Range("S6").Value = 0.2683720948477344
$ws.Range("T6").Value = 0.2683720948477344

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntn1"
$ws.Range("C7").Value = "Unc5a"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 44.154177
$ws.Range("H7").Value = 132.462531
$ws.Range("I7").Value = 0.8392647337471152
$ws.Range("J7").Value = 0.8392647337471153
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8581799999999999
$ws.Range("N7").Value = 2.57454
$ws.Range("O7").Value = 0.18526016866077
$ws.Range("P7").Value = 0.1852601686607701
$ws.Range("Q7").Value = 37.89223161786
$ws.Range("R7").Value = 341.03008456074
$ws.Range("S7").Value = 0.1554823261250268
$ws.Range("T7").Value = 0.1554823261250269

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ntn1"
$ws.Range("C8").Value = "Unc5a"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.525807
$ws.Range("H8").Value = 22.577421
$ws.Range("I8").Value = 0.1430474948743168
$ws.Range("J8").Value = 0.1430474948743168
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.292844666666667
$ws.Range("N8").Value = 6.878534
$ws.Range("O8").Value = 0.4949693416994264
$ws.Range("P8").Value = 0.4949693416994265
$ws.Range("Q8").Value = 17.25550644231267
$ws.Range("R8").Value = 155.299557980814
$ws.Range("S8").Value = 0.07080412436969265
$ws.Range("T8").Value = 0.07080412436969266

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ntn1"
$ws.Range("C9").Value = "Unc5a"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.525807
$ws.Range("H9").Value = 22.577421
$ws.Range("I9").Value = 0.1430474948743168
$ws.Range("J9").Value = 0.1430474948743168
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.481271666666667
$ws.Range("N9").Value = 4.443815
$ws.Range("O9").Value = 0.3197704896398035
$ws.Range("P9").Value = 0.3197704896398036
$ws.Range("Q9").Value = 11.14776467790167
$ws.Range("R9").Value = 100.329882101115
$ws.Range("S9").Value = 0.04574236747770757
$ws.Range("T9").Value = 0.04574236747770757

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ntn1"
$ws.Range("C10").Value = "Unc5a"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.525807
$ws.Range("H10").Value = 22.577421
$ws.Range("I10").Value = 0.1430474948743168
$ws.Range("J10").Value = 0.1430474948743168
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8581799999999999
$ws.Range("N10").Value = 2.57454
$ws.Range("O10").Value = 0.18526016866077
$ws.Range("P10").Value = 0.1852601686607701
$ws.Range("Q10").Value = 6.45849705126
$ws.Range("R10").Value = 58.12647346134
$ws.Range("S10").Value = 0.02650100302691656
$ws.Range("T10").Value = 0.02650100302691657

Write-Output "Edit complete"